$wb = $excel.ActiveWorkbook

# --- Sheet1: append leaderboard rows 22-33 ---
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Cells.Item(22, 1).Value = "Jack"
$ws1.Cells.Item(22, 2).Value = 521

$ws1.Cells.Item(23, 1).Value = "Jack"
$ws1.Cells.Item(23, 2).Value = 0

$ws1.Cells.Item(24, 1).Value = "Jack"
$ws1.Cells.Item(24, 2).Value = 1563

$ws1.Cells.Item(25, 1).Value = "p"
$ws1.Cells.Item(25, 2).Value = 1042

$ws1.Cells.Item(26, 1).Value = "a"
$ws1.Cells.Item(26, 2).Value = 4168

$ws1.Cells.Item(27, 1).Value = "q"
$ws1.Cells.Item(27, 2).Value = 1042

$ws1.Cells.Item(28, 1).Value = "a"
$ws1.Cells.Item(28, 2).Value = 1042

$ws1.Cells.Item(29, 1).Value = "a"
$ws1.Cells.Item(29, 2).Value = 1042

$ws1.Cells.Item(30, 1).Value = "j"
$ws1.Cells.Item(30, 2).Value = 1042

$ws1.Cells.Item(31, 1).Value = "l"
$ws1.Cells.Item(31, 2).Value = 1042

$ws1.Cells.Item(32, 1).Value = "a"
$ws1.Cells.Item(32, 2).Value = 1042

$ws1.Cells.Item(33, 1).Value = "a"
$ws1.Cells.Item(33, 2).Value = 0

# --- Sheet4: append leaderboard rows 7-9 ---
$ws4 = $wb.Worksheets.Item("Sheet4")

$ws4.Cells.Item(7, 1).Value = "l"
$ws4.Cells.Item(7, 2).Value = 1042

$ws4.Cells.Item(8, 1).Value = "a"
$ws4.Cells.Item(8, 2).Value = 0

$ws4.Cells.Item(9, 1).Value = "Jack"
$ws4.Cells.Item(9, 2).Value = 521

# --- Sheet5: append leaderboard row 7 ---
$ws5 = $wb.Worksheets.Item("Sheet5")

$ws5.Cells.Item(7, 1).Value = "Jack"
$ws5.Cells.Item(7, 2).Value = 16151
